# Insert a new pricing record for "Vega Modelo de Temuco" (Uva) as row 466,
# pushing all existing rows from 466 downward by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 466 (shifts 466..527 down to 467..528).
$ws.Rows.Item(466).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(466, 1).Value = 10
$ws.Cells.Item(466, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(466, 3).Value = "La Araucanía"
$ws.Cells.Item(466, 4).Value = 44474
$ws.Cells.Item(466, 5).Value = 9
$ws.Cells.Item(466, 6).Value = "Fruta"
$ws.Cells.Item(466, 7).Value = 100109
$ws.Cells.Item(466, 8).Value = "Uva"
$ws.Cells.Item(466, 9).Value = 100109001
$ws.Cells.Item(466, 10).Value = "Uva"
$ws.Cells.Item(466, 11).Value = "Superior Seedless"
$ws.Cells.Item(466, 12).Value = "Primera"
$ws.Cells.Item(466, 13).Value = 300
$ws.Cells.Item(466, 14).Value = 37000
$ws.Cells.Item(466, 15).Value = 37000
$ws.Cells.Item(466, 16).Value = 37000
$ws.Cells.Item(466, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(466, 18).Value = "EE.UU."
$ws.Cells.Item(466, 19).Value = 4625
$ws.Cells.Item(466, 20).Value = 8
